# print_request_logsheet.xlsx - "July 2018" sheet
#
# - Mark the two existing "Vaulted IR Phantom" mold pieces (rows 2 & 3) as
#   completed on 03-07-2018.
# - Log a new print request on row 4 for an "EVHP Concept 5 TEE Track Wide
#   Bore", requested 04-07-2018.
#
# Note: the date strings in this sheet ("dd-mm-yyyy") are stored as plain
# text, not real Excel dates - typing them straight into a cell via
# .Value would normally get auto-parsed into a serial date number. To
# keep them as literal text (matching the rest of the sheet) we build the
# string on a scratch cell as a formula result (a plain string, never
# date-sniffed) and bring only the *value* across with PasteSpecial, then
# bring only the *format* across separately from an existing, correctly
# styled cell. This also avoids minting any new cell styles/number
# formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

$xlPasteValues = -4163
$xlPasteFormats = -4122

$scratch = $ws.Range("AZ1000")

function Set-LiteralText {
    param($TargetCell, [string]$Text)

    # Writing ="<text>" always yields a plain string result -> never
    # reinterpreted as a date/number, regardless of its shape.
    $scratch.Value = '="' + $Text + '"'
    $scratch.Copy()
    $TargetCell.PasteSpecial($xlPasteValues)
    # Full Clear (not just ClearContents) so the scratch cell leaves no
    # trace (value *or* inherited column style) in the saved sheet.
    $scratch.Clear()
}

# --- New request: row 4 (date entered first -> lands lower in the shared
#     string table, matching the order the workbook was edited in) ---
Set-LiteralText $ws.Range("A4") "04-07-2018"
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)

# --- Mark rows 2 & 3 ("Vaulted IR Phantom - Body" / "- Ends") complete ---
Set-LiteralText $ws.Range("B2") "03-07-2018"
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial($xlPasteFormats)

Set-LiteralText $ws.Range("B3") "03-07-2018"
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial($xlPasteFormats)

$ws.Range("C4").Value = "EVHP Concept 5 TEE Track Wide Bore"
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial($xlPasteFormats)

$ws.Range("D4").Value = 1
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial($xlPasteFormats)

$ws.Range("E4").Value = "Polylite"
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial($xlPasteFormats)

$ws.Range("F4").Value = 2
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial($xlPasteFormats)

$ws.Range("G4").Value = 20
$ws.Range("G2").Copy()
$ws.Range("G4").PasteSpecial($xlPasteFormats)

$ws.Range("H4").Value = 0.2
$ws.Range("H2").Copy()
$ws.Range("H4").PasteSpecial($xlPasteFormats)

$ws.Range("I4").Value = "NA"
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)

# Leave the cursor where data entry finished, as in the saved workbook.
$ws.Range("I4").Select() | Out-Null
